$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for numeric-looking strings in D:E, matching the
# original inlineStr cell type (Excel would otherwise auto-convert strings
# like "1.009" into numbers). Style is restored to Normal afterwards so no
# extra style index is left attached to the data cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "28.512.07"
$ws.Cells.Item(2, 5).Value = "  -3.51%  "
$ws.Cells.Item(3, 4).Value = "1.958.48"
$ws.Cells.Item(3, 5).Value = "  -2.23%  "
$ws.Cells.Item(4, 5).Value = "  -0.36%  "
$ws.Cells.Item(5, 4).Value = "321.29"
$ws.Cells.Item(5, 5).Value = "  -2.47%  "
$ws.Cells.Item(6, 4).Value = "1.009"
$ws.Cells.Item(6, 5).Value = "  -0.32%  "
$ws.Cells.Item(7, 4).Value = "0.4766"
$ws.Cells.Item(7, 5).Value = "  -4.69%  "
$ws.Cells.Item(8, 4).Value = "0.4070"
$ws.Cells.Item(8, 5).Value = "  -3.52%  "
$ws.Cells.Item(9, 4).Value = "53.51"
$ws.Cells.Item(9, 5).Value = "  -1.84%  "
$ws.Cells.Item(10, 4).Value = "0.08450"
$ws.Cells.Item(10, 5).Value = "  -6.33%  "
$ws.Cells.Item(11, 4).Value = "1.061"
$ws.Cells.Item(11, 5).Value = "  -5.03%  "
$ws.Cells.Item(12, 4).Value = "22.13"
$ws.Cells.Item(12, 5).Value = "  -4.90%  "
$ws.Cells.Item(13, 4).Value = "1.981.97"
$ws.Cells.Item(13, 5).Value = "  -3.73%  "
$ws.Cells.Item(14, 4).Value = "7.636"
$ws.Cells.Item(14, 5).Value = "  -5.02%  "
$ws.Cells.Item(15, 4).Value = "6.193"
$ws.Cells.Item(15, 5).Value = "  -4.22%  "
$ws.Cells.Item(16, 4).Value = "1.011"
$ws.Cells.Item(16, 5).Value = "  -0.30%  "
$ws.Cells.Item(17, 4).Value = "0.00001074"
$ws.Cells.Item(17, 5).Value = "  -3.65%  "
$ws.Cells.Item(18, 4).Value = "89.33"
$ws.Cells.Item(18, 5).Value = "  -5.40%  "
$ws.Cells.Item(19, 4).Value = "0.06626"
$ws.Cells.Item(19, 5).Value = "  -0.95%  "
$ws.Cells.Item(20, 4).Value = "18.76"
$ws.Cells.Item(20, 5).Value = "  -4.56%  "
$ws.Cells.Item(21, 5).Value = "  -0.21%  "
$ws.Cells.Item(22, 4).Value = "5.822"
$ws.Cells.Item(22, 5).Value = "  -2.57%  "
$ws.Cells.Item(23, 4).Value = "28.514.13"
$ws.Cells.Item(23, 5).Value = "  -3.66%  "
$ws.Cells.Item(24, 4).Value = "11.57"
$ws.Cells.Item(24, 5).Value = "  -3.46%  "
$ws.Cells.Item(25, 4).Value = "2.289"
$ws.Cells.Item(25, 5).Value = "  -0.62%  "
$ws.Cells.Item(26, 4).Value = "2.200.62"
$ws.Cells.Item(26, 5).Value = "  -4.11%  "
$ws.Cells.Item(27, 4).Value = "154.30"
$ws.Cells.Item(27, 5).Value = "  -2.81%  "
$ws.Cells.Item(28, 4).Value = "20.24"
$ws.Cells.Item(28, 5).Value = "  -2.31%  "
$ws.Cells.Item(29, 4).Value = "5.998"
$ws.Cells.Item(29, 5).Value = "  -5.28%  "
$ws.Cells.Item(30, 4).Value = "2.166"
$ws.Cells.Item(30, 5).Value = "  -5.63%  "
$ws.Cells.Item(31, 4).Value = "123.96"
$ws.Cells.Item(31, 5).Value = "  -3.33%  "
$ws.Cells.Item(32, 4).Value = "0.9918"
$ws.Cells.Item(32, 5).Value = "  -5.94%  "
$ws.Cells.Item(33, 4).Value = "0.09602"
$ws.Cells.Item(33, 5).Value = "  -3.53%  "
$ws.Cells.Item(34, 4).Value = "1.455"
$ws.Cells.Item(34, 5).Value = "  -7.10%  "
$ws.Cells.Item(35, 4).Value = "5.600"
$ws.Cells.Item(35, 5).Value = "  -3.96%  "
$ws.Cells.Item(36, 4).Value = "3.651"
$ws.Cells.Item(36, 5).Value = "  -3.93%  "
$ws.Cells.Item(37, 4).Value = "0.02334"
$ws.Cells.Item(37, 5).Value = "  -5.34%  "
$ws.Cells.Item(38, 4).Value = "8.812"
$ws.Cells.Item(38, 5).Value = "  -4.93%  "
$ws.Cells.Item(39, 4).Value = "0.06229"
$ws.Cells.Item(39, 5).Value = "  -2.72%  "
$ws.Cells.Item(40, 4).Value = "1.258"
$ws.Cells.Item(40, 5).Value = "  -3.55%  "
$ws.Cells.Item(41, 4).Value = "0.6236"
$ws.Cells.Item(41, 5).Value = "  -4.63%  "
$ws.Cells.Item(42, 4).Value = "11.16"
$ws.Cells.Item(42, 5).Value = "  -4.35%  "
$ws.Cells.Item(43, 4).Value = "1.009"
$ws.Cells.Item(43, 5).Value = "  -0.32%  "
$ws.Cells.Item(44, 4).Value = "0.1927"
$ws.Cells.Item(44, 5).Value = "  -5.80%  "
$ws.Cells.Item(45, 4).Value = "1.335"
$ws.Cells.Item(45, 5).Value = "  +2.26%  "
$ws.Cells.Item(46, 4).Value = "0.5976"
$ws.Cells.Item(46, 5).Value = "  -5.80%  "
$ws.Cells.Item(47, 4).Value = "13.06"
$ws.Cells.Item(47, 5).Value = "  -2.84%  "
$ws.Cells.Item(48, 4).Value = "2.060"
$ws.Cells.Item(48, 5).Value = "  -6.03%  "
$ws.Cells.Item(49, 4).Value = "3.401"
$ws.Cells.Item(49, 5).Value = "  -3.06%  "
$ws.Cells.Item(50, 4).Value = "0.00000000329"
$ws.Cells.Item(50, 5).Value = "  -1.84%  "
$ws.Cells.Item(51, 4).Value = "0.06840"
$ws.Cells.Item(51, 5).Value = "  -2.13%  "

$ws.Range("D2:E51").Style = "Normal"
